# Apply the "assembler funcionando" edit to commands.xlsx
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3 header ("Integer Arithmetic Operations") restyle: A3:E3 ---
# New style (cellXfs idx 26) reuses font 9 (bold 16 PetitaLight, RGB 000000)
# instead of font 8 (bold 16 PetitaLight, theme color 1), same center/center
# alignment. Setting Font.Color to pure black on the already-centered range
# makes the engine materialize that exact new xf.
$ws.Range("A3:E3").Font.Color = 0

# --- New "Operation" sub-header row above the System/Program tables (A46) ---
# Clone A3's freshly-built style onto A46:E46 (format-only paste) so it
# reuses the same cellXfs entry instead of growing a second, slightly
# different one, then drop in the label text and merge.
$ws.Range("A3").Copy()
$ws.Range("A46:E46").PasteSpecial(-4122)
$ws.Range("A46").Value = "      Operation"
$ws.Range("A46:E46").Merge()
$ws.Rows.Item(46).RowHeight = 27

# --- Assembler column ("F") progress marks: no -> yes ---
$ws.Range("F4").Value = "yes`nyes`nyes"
$ws.Range("F5").Value = "yes`nyes`nyes"
$ws.Range("F6").Value = "yes`nyes"

# Row 7: collapse the stray double "no`nno" down to a single "no" in F/G/H,
# matching the single-line entries used by every other row, and shrink the
# row back to the single-line height.
$ws.Range("F7:H7").Value = "no"
$ws.Rows.Item(7).RowHeight = 28

$ws.Range("F9").Value = "yes"

# --- Selection left where the user's cursor ended up ---
$ws.Range("D4").Select()
